$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "Administrator, Miss Dina Nasr"
$replacement = "Miss Dina Nasr, Administrator"

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $target) {
        $cell.Value2 = $replacement
    }
}
